# "filter for trademark cost calc"
# - Inserts a new data row (for address 0xa4ada63fa9417e00c2794cd51db7e9a6aaa473d9)
#   before the current row 42, pushing the existing rows 42-49 down to 43-50.
# - Bumps four pre-existing "mint" counter cells by 1 (N3, S3, X3, N24) as a
#   side effect of the new address being folded into those aggregates.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. insert a blank row at row 42, shifting rows 42:49 down to 43:50 ---
$ws.Rows.Item(42).Insert()

# --- 2. populate the newly inserted row 42 ---
# All cells in this sheet are plain text (even the numeric-looking ones), so
# every value is written with a leading apostrophe to force text storage
# instead of Excel's automatic number inference.
$newRowValues = @(
    "0xa4ada63fa9417e00c2794cd51db7e9a6aaa473d9",
    "3",
    "0",
    "30",
    "0",
    "27",
    "0.00",
    "29",
    "0",
    "0",
    "29",
    "3.11",
    "29",
    "0",
    "0",
    "29",
    "3.19",
    "29",
    "1",
    "0",
    "29",
    "0.10",
    "31",
    "0",
    "0",
    "29",
    "0.08"
)

for ($i = 0; $i -lt $newRowValues.Length; $i++) {
    $ws.Cells.Item(42, $i + 1).Value = "'" + $newRowValues[$i]
}

# --- 3. bump the four scattered mint counters by 1 ---
$ws.Range("N3").Value  = "'38"
$ws.Range("S3").Value  = "'90"
$ws.Range("X3").Value  = "'98"
$ws.Range("N24").Value = "'4"
